$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $cellAddr, $val) {
    $sheet.Range($cellAddr).NumberFormat = "@"
    $sheet.Range($cellAddr).Value = $val
    $sheet.Range($cellAddr).Style = "Normal"
}

$ws.Range("D2").Value = "37.775.46"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "2.083.67"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue $ws "D5" "233.72"
$ws.Range("E5").Value = "  -0.98%  "
Set-TextValue $ws "D6" "0.624"
$ws.Range("E6").Value = "  -0.62%  "
Set-TextValue $ws "D7" "58.76"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.12%  "
Set-TextValue $ws "D10" "0.0787"
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("D12").Value = "2.388.51"
$ws.Range("E12").Value = "  -1.83%  "
Set-TextValue $ws "D13" "14.86"
$ws.Range("E13").Value = "  +1.36%  "
Set-TextValue $ws "D14" "21.04"
$ws.Range("E14").Value = "  -3.40%  "
Set-TextValue $ws "D15" "0.777"
$ws.Range("E15").Value = "  -2.35%  "
Set-TextValue $ws "D16" "5.37"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "2.092.77"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "37.729.02"
$ws.Range("E18").Value = "  -1.34%  "
Set-TextValue $ws "D19" "6.13"
$ws.Range("E19").Value = "  -0.57%  "
Set-TextValue $ws "D20" "71.68"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").Value = "  +0.69%  "
Set-TextValue $ws "D22" "228.59"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -0.50%  "
Set-TextValue $ws "D25" "2.34"
$ws.Range("E25").Value = "  -3.12%  "
Set-TextValue $ws "D26" "171.49"
$ws.Range("E26").Value = "  +1.63%  "
Set-TextValue $ws "D27" "9.13"
$ws.Range("E27").Value = "  +0.49%  "
Set-TextValue $ws "D28" "0.136"
$ws.Range("E28").Value = "  -3.68%  "
Set-TextValue $ws "D29" "19.54"
$ws.Range("E29").Value = "  -0.39%  "
Set-TextValue $ws "D30" "1.41"
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("E31").Value = "  +1.58%  "
Set-TextValue $ws "D32" "4.70"
$ws.Range("E32").Value = "  +0.62%  "
Set-TextValue $ws "D33" "0.0635"
$ws.Range("E33").Value = "  +0.82%  "
Set-TextValue $ws "D34" "4.69"
$ws.Range("E34").Value = "  +1.27%  "
Set-TextValue $ws "D35" "2.48"
$ws.Range("E35").Value = "  -5.23%  "
Set-TextValue $ws "D36" "1.82"
$ws.Range("E36").Value = "  -0.50%  "
Set-TextValue $ws "D37" "3.42"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("E38").Value = "  -0.11%  "
Set-TextValue $ws "D39" "5.37"
$ws.Range("E39").Value = "  -3.19%  "
$ws.Range("B40").Value = "Cronos"
$ws.Range("C40").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D40" "0.0977"
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D41" "99.93"
$ws.Range("E41").Value = "  +2.30%  "
$ws.Range("E42").Value = "  -2.64%  "
Set-TextValue $ws "D43" "0.0216"
$ws.Range("E43").Value = "  -0.18%  "
Set-TextValue $ws "D44" "16.80"
$ws.Range("E44").Value = "  +3.23%  "
$ws.Range("D45").Value = "1.441.88"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("E46").Value = "  -1.22%  "
Set-TextValue $ws "D47" "4.19"
$ws.Range("E47").Value = "  +0.04%  "
Set-TextValue $ws "D48" "1.07"
$ws.Range("E48").Value = "  -1.31%  "
Set-TextValue $ws "D49" "7.42"
$ws.Range("E49").Value = "  +1.37%  "
Set-TextValue $ws "D50" "3.00"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Value = "2.273.10"
$ws.Range("E51").Value = "  -1.91%  "